# Adds a new "intervention_type" column (K) to the clinical-trials list,
# one DRUG/DEVICE/... value per trial row (2..121). Mirrors the commit's
# addition of an "intervention_type" indicator for every sponsor row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cell K1: same label style as the other header cells (A1:J1) ---
$ws.Range("K1").Value = "intervention_type"
$ws.Range("A1").Copy()
$ws.Range("K1").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# --- Data rows K2:K121, in row order ---
$interventionTypes = @(
    'DRUG','DRUG','DRUG','DRUG','DEVICE','DEVICE','DRUG','DRUG','DEVICE','GENETIC',
    'DRUG','BEHAVIORAL','OTHER','PROCEDURE','','DRUG','DRUG','','BEHAVIORAL','DEVICE',
    'OTHER','DEVICE','RADIATION','DRUG','DEVICE','DRUG','DEVICE','PROCEDURE','OTHER','PROCEDURE',
    'DRUG','DRUG','DEVICE','DEVICE','PROCEDURE','DRUG','DIAGNOSTIC_TEST','DEVICE','OTHER','DEVICE',
    'DEVICE','DRUG','DRUG','RADIATION','DEVICE','BIOLOGICAL','PROCEDURE','DEVICE','OTHER','DEVICE',
    'DEVICE','OTHER','DEVICE','DRUG','BEHAVIORAL','DRUG','PROCEDURE','DIAGNOSTIC_TEST','DIAGNOSTIC_TEST','DRUG',
    'BEHAVIORAL','DEVICE','DRUG','OTHER','PROCEDURE','DIAGNOSTIC_TEST','OTHER','PROCEDURE','OTHER','DRUG',
    'RADIATION','DIAGNOSTIC_TEST','DEVICE','OTHER','PROCEDURE','DEVICE','DEVICE','DIAGNOSTIC_TEST','OTHER','DIAGNOSTIC_TEST',
    'DEVICE','DRUG','DIAGNOSTIC_TEST','DIAGNOSTIC_TEST','DEVICE','OTHER','DEVICE','OTHER','','DEVICE',
    'PROCEDURE','OTHER','OTHER','PROCEDURE','OTHER','DEVICE','DEVICE','DEVICE','DIAGNOSTIC_TEST','DEVICE',
    'DRUG','PROCEDURE','DRUG','DRUG','DRUG','DRUG','OTHER','DRUG','DRUG','DRUG',
    'PROCEDURE','DRUG','DRUG','OTHER','OTHER','DRUG','PROCEDURE','DRUG','DRUG','PROCEDURE'
)

for ($i = 0; $i -lt $interventionTypes.Length; $i++) {
    $r = $i + 2
    $val = $interventionTypes[$i]
    # Rows 16, 19 and 90 have no known intervention type in the source data
    # (left blank), matching every other column's "unknown" convention.
    if ($val -ne "") {
        $ws.Cells.Item($r, 11).Value = $val
    }
}
